$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that is refreshed on
# every automatic update run. Bump it from 45189 (2023-09-20) to
# 45190 (2023-09-21) for every data row (rows 2 through 339).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 339 }

$ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item($lastRow, 3)).Value = 45190
